# Apply the "0.1 -> 1.0" release edits described in the commit message.
$d = $word.ActiveDocument

# 1. Version number in the "Changes" history table: 0.1 -> 1.0
$d.Content.Find.Execute("0.1", $true, $false, $false, $false, $false, $true, 1, $false, "1.0", 2) | Out-Null

# 2. Change description in the same table row: Criacao -> Alteracao
$d.Content.Find.Execute("Criacao", $true, $false, $false, $false, $false, $true, 1, $false, "Alteracao", 2) | Out-Null

# 3/4. Lower-case "permissoes" -> "Permissoes" in the basic flow text and
#      in the postconditions table (both occurrences of the exact word,
#      replaced in a single wdReplaceAll pass).
$d.Content.Find.Execute("permissoes", $true, $false, $false, $false, $false, $true, 1, $false, "Permissoes", 2) | Out-Null
